$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet,
#    copying the row1/row2 formatting from "2021-Q3" so borders/fonts
#    match the other quarterly sheets.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1_2022 = $wb.Worksheets.Add($totalSheet)
$q1_2022.Name = "2022-Q1"

$q3_2021 = $wb.Worksheets.Item("2021-Q3")
$q3_2021.Range("A1:H2").Copy()
$q1_2022.Range("A1:H2").PasteSpecial(-4122)

# Header row
$q1_2022.Range("B1").Value = "基金代码"
$q1_2022.Range("C1").Value = "基金名称"
$q1_2022.Range("D1").Value = "基金规模"
$q1_2022.Range("E1").Value = "股票总仓位"
$q1_2022.Range("F1").Value = "仓位占比"
$q1_2022.Range("G1").Value = "持有市值(亿元)"
$q1_2022.Range("H1").Value = "仓位排名"

# Data row 2 - text-like numeric fields need to stay text, so mark the
# range as text before assigning them (otherwise Excel silently
# re-interprets "0.14" etc. as numbers).
$q1_2022.Range("A2").Value = 0
$q1_2022.Range("B2:G2").NumberFormat = "@"
$q1_2022.Range("B2").Value = "539002"
$q1_2022.Range("C2").Value = "建信新兴市场优选混合QDII"
$q1_2022.Range("D2").Value = "0.14"
$q1_2022.Range("E2").Value = "83.76"
$q1_2022.Range("F2").Value = "5.07"
$q1_2022.Range("G2").Value = "0.0071"
$q1_2022.Range("H2").Value = 5

# Restore the originally-active sheet/tab.
$wb.Worksheets.Item("2021-Q1").Activate()

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new first data row for 2022-Q1
#    and shift the existing quarters down by one row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Grow the used range by one row, copying row 4's formatting down into
# the new row 5 so the new row's "index" cell (A5) gets the same
# border/bold/centered style as A2:A4.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

# Shift the existing 3 data rows down by one (bottom-up so a row isn't
# overwritten before it has been read).
$total.Range("A5").Value = 3
$total.Range("B5").Value = $total.Range("B4").Value()
$total.Range("C5").Value = $total.Range("C4").Value()
$total.Range("D5").Value = $total.Range("D4").Value()

$total.Range("A4").Value = 2
$total.Range("B4").Value = $total.Range("B3").Value()
$total.Range("C4").Value = $total.Range("C3").Value()
$total.Range("D4").Value = $total.Range("D3").Value()

$total.Range("A3").Value = 1
$total.Range("B3").Value = $total.Range("B2").Value()
$total.Range("C3").Value = $total.Range("C2").Value()
$total.Range("D3").Value = $total.Range("D2").Value()

# Write the brand-new 2022-Q1 row into row 2.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01
